$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 56.88889
$ws.Range("I5").Value = 57.75
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 57.75
$ws.Range("L5").Value = 50
$ws.Range("M5").Value = 57.25
$ws.Range("N5").Value = -280

$ws.Range("H9").Value = 250229.75
$ws.Range("I9").Value = 500074.5
$ws.Range("J9").Value = 385
$ws.Range("K9").Value = 500074.5
$ws.Range("L9").Value = 385
$ws.Range("M9").Value = -499905.5

$ws.Range("H19").Value = 701.6667
$ws.Range("I19").Value = 1000
$ws.Range("J19").Value = 552.5
$ws.Range("K19").Value = 1000
$ws.Range("L19").Value = 552.5
$ws.Range("M19").Value = -825
$ws.Range("N19").Value = -902.5

$ws.Range("H28").Value = 410.8
$ws.Range("I28").Value = 410.8
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 410.8
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 74.19999999999999

$ws.Range("H53").Value = 302.2857
$ws.Range("I53").Value = 358.66666
$ws.Range("J53").Value = 260
$ws.Range("K53").Value = 358.66666
$ws.Range("L53").Value = 260
$ws.Range("M53").Value = 278.33334
$ws.Range("N53").Value = -1534

$ws.Range("H98").Value = 2481.1667
$ws.Range("I98").Value = 749.875
$ws.Range("J98").Value = 5943.75
$ws.Range("K98").Value = 749.875
$ws.Range("L98").Value = 5943.75
$ws.Range("M98").Value = 748.125

$ws.Range("H100").Value = 3208.5334
$ws.Range("I100").Value = 3366.2856
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 3366.2856
$ws.Range("L100").Value = 1000
$ws.Range("M100").Value = -2825.2856

$ws.Range("H107").Value = 1464.421
$ws.Range("I107").Value = 1434.6666
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 1434.6666
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = 485.3334

$ws.Range("H122").Value = 2481.1667
$ws.Range("I122").Value = 749.875
$ws.Range("J122").Value = 5943.75
$ws.Range("K122").Value = 2249.625
$ws.Range("L122").Value = 17831.25
$ws.Range("M122").Value = 200.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3420.3333
$ws.Range("I61").Value = 3420.3333
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3420.3333
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3208.3333

$ws.Range("M132").ClearContents()
$ws.Range("H132").Value = 4224.25
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4224.25
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 12672.75
$ws.Range("N132").Value = -17732.75

$ws.Range("H136").Value = 3420.3333
$ws.Range("I136").Value = 3420.3333
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10260.9999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -7710.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1073.5
$ws.Range("I134").Value = 1073.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3220.5
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -685.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2207.1538
$ws.Range("I31").Value = 1445.8572
$ws.Range("J31").Value = 3095.3333
$ws.Range("K31").Value = 1445.8572
$ws.Range("L31").Value = 3095.3333
$ws.Range("M31").Value = -1150.8572
$ws.Range("N31").Value = -3685.3333

$ws.Range("H34").Value = 2207.1538
$ws.Range("I34").Value = 1445.8572
$ws.Range("J34").Value = 3095.3333
$ws.Range("K34").Value = 1445.8572
$ws.Range("L34").Value = 3095.3333
$ws.Range("M34").Value = -1243.8572
$ws.Range("N34").Value = -3499.3333

$ws.Range("H132").Value = 3797.4285
$ws.Range("I132").Value = 3296
$ws.Range("J132").Value = 4173.5
$ws.Range("K132").Value = 9888
$ws.Range("L132").Value = 12520.5
$ws.Range("M132").Value = -7358

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 6599.8
$ws.Range("I34").Value = 499.5
$ws.Range("J34").Value = 10666.667
$ws.Range("K34").Value = 1498.5
$ws.Range("L34").Value = 32000.001
$ws.Range("M34").Value = -1414.5
$ws.Range("N34").Value = -32168.001

$ws.Range("H39").Value = 19495
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 19495
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 58485
$ws.Range("N39").Value = -59073

$ws.Range("N46").ClearContents()
$ws.Range("H46").Value = 100
$ws.Range("I46").Value = 100
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 300
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -209

$ws.Range("M50").ClearContents()
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0

$ws.Range("M53").ClearContents()
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0

$ws.Range("H55").Value = 13997.375
$ws.Range("I55").Value = 5000
$ws.Range("J55").Value = 15282.714
$ws.Range("K55").Value = 15000
$ws.Range("L55").Value = 45848.142
$ws.Range("M55").Value = -14823
$ws.Range("N55").Value = -46202.142

$ws.Range("N58").ClearContents()
$ws.Range("H58").Value = 1850
$ws.Range("I58").Value = 1850
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 5550
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -5422

$ws.Range("H82").Value = 15000
$ws.Range("I82").Value = 15000
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 45000
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -44594

$ws.Range("H85").Value = 15000
$ws.Range("I85").Value = 15000
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 45000
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -43596

$ws.Range("H94").Value = 3750
$ws.Range("I94").Value = 1666.6666
$ws.Range("J94").Value = 10000
$ws.Range("K94").Value = 4999.9998
$ws.Range("L94").Value = 30000
$ws.Range("M94").Value = -4323.9998
$ws.Range("N94").Value = -31352

$ws.Range("N103").ClearContents()
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 0

$ws.Range("H106").Value = 5000
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 5000
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 15000
$ws.Range("N106").Value = -16892

$ws.Range("H109").Value = 4351.607
$ws.Range("I109").Value = 1974.5
$ws.Range("J109").Value = 4999.909
$ws.Range("K109").Value = 5923.5
$ws.Range("L109").Value = 14999.727
$ws.Range("M109").Value = -4883.5
$ws.Range("N109").Value = -17079.727

$ws.Range("M112").ClearContents()
$ws.Range("H112").Value = 7500
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 7500
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 22500
$ws.Range("N112").Value = -24716

$ws.Range("H136").Value = 2370.3125
$ws.Range("I136").Value = 2012.5
$ws.Range("J136").Value = 2421.4285
$ws.Range("K136").Value = 6037.5
$ws.Range("L136").Value = 7264.2855
$ws.Range("M136").Value = -937.5
$ws.Range("N136").Value = -17464.2855

$ws.Range("H138").Value = 3759.75
$ws.Range("I138").Value = 2245
$ws.Range("J138").Value = 5274.5
$ws.Range("K138").Value = 6735
$ws.Range("L138").Value = 15823.5
$ws.Range("M138").Value = -1595
$ws.Range("N138").Value = -26103.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 800
$ws.Range("I43").Value = 800
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 800
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -649

$ws.Range("N101").ClearContents()
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 0

$ws.Range("H107").Value = 6100
$ws.Range("I107").Value = 7500
$ws.Range("J107").Value = 3300
$ws.Range("K107").Value = 7500
$ws.Range("L107").Value = 3300
$ws.Range("M107").Value = -5580
$ws.Range("N107").Value = -7140

$ws.Range("H126").Value = 10000
$ws.Range("I126").Value = 10000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 30000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -27530

$ws.Range("H132").Value = 4475.5454
$ws.Range("I132").Value = 3258.8
$ws.Range("J132").Value = 5489.5
$ws.Range("K132").Value = 9776.400000000001
$ws.Range("L132").Value = 16468.5
$ws.Range("M132").Value = -7246.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3638.3333
$ws.Range("I7").Value = 3638.3333
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 3638.3333
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -3526.3333

$ws.Range("H22").Value = 4813.6113
$ws.Range("I22").Value = 2687
$ws.Range("J22").Value = 7471.875
$ws.Range("K22").Value = 2687
$ws.Range("L22").Value = 7471.875
$ws.Range("M22").Value = -2392
$ws.Range("N22").Value = -8061.875

$ws.Range("H27").Value = 4813.6113
$ws.Range("I27").Value = 2687
$ws.Range("J27").Value = 7471.875
$ws.Range("K27").Value = 2687
$ws.Range("L27").Value = 7471.875
$ws.Range("M27").Value = -2580
$ws.Range("N27").Value = -7685.875

$ws.Range("M46").ClearContents()
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0

$ws.Range("H100").Value = 1449.75
$ws.Range("I100").Value = 1366.3334
$ws.Range("J100").Value = 1700
$ws.Range("K100").Value = 1366.3334
$ws.Range("L100").Value = 1700
$ws.Range("M100").Value = -825.3334
$ws.Range("N100").Value = -2782

$ws.Range("H101").Value = 15000
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 15000
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 15000
$ws.Range("N101").Value = -21490

$ws.Range("H126").Value = 3638.3333
$ws.Range("I126").Value = 3638.3333
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 10914.9999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -8444.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1629.4546
$ws.Range("I126").Value = 1469.6
$ws.Range("J126").Value = 1676.4706
$ws.Range("K126").Value = 4408.799999999999
$ws.Range("L126").Value = 5029.4118
$ws.Range("M126").Value = -1938.799999999999
$ws.Range("N126").Value = -9969.4118

$ws.Range("H132").Value = 3204.9
$ws.Range("I132").Value = 1626.6666
$ws.Range("J132").Value = 7939.6
$ws.Range("K132").Value = 4879.9998
$ws.Range("L132").Value = 23818.8
$ws.Range("M132").Value = -2349.9998

$ws.Range("H136").Value = 852.8889
$ws.Range("I136").Value = 848.5
$ws.Range("J136").Value = 888
$ws.Range("K136").Value = 2545.5
$ws.Range("L136").Value = 2664
$ws.Range("M136").Value = 4.5
